$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1205, 477, "V", "NewOrleans", 0, 10, 3, 3, 16, 42.5, 40.5, 160, 20.5),
  @(1205, 478, "H", "TampaBay", 3, 0, 0, 14, 17, 7, 3, -180, 3.5),
  @(1208, 101, "V", "LasVegas", 10, 3, 0, 3, 16, 4.5, 7, -260, 20.5),
  @(1208, 102, "H", "LARams", 3, 0, 0, 14, 17, 43, 42, 220, 0.5),
  @(1211, 105, "V", "Minnesota", 7, 0, 6, 10, 23, 3, 51.5, 120, 3),
  @(1211, 106, "H", "Detroit", 7, 7, 7, 13, 34, 51, 1, -140, 26),
  @(1211, 107, "V", "NYJets", 0, 0, 7, 5, 12, 44.5, 43.5, 375, 20.5),
  @(1211, 108, "H", "Buffalo", 0, 7, 10, 3, 20, 9.5, 10.5, -450, 3),
  @(1211, 109, "V", "Baltimore", 10, 3, 0, 3, 16, 4, 36.5, 105, 17.5),
  @(1211, 110, "H", "Pittsburgh", 7, 0, 0, 7, 14, 41, 1, -125, 1),
  @(1211, 111, "V", "Philadelphia", 7, 17, 10, 14, 48, 5.5, 8.5, -330, 2.5),
  @(1211, 112, "H", "NYGiants", 0, 7, 7, 8, 22, 43, 44, 270, 20.5),
  @(1211, 113, "V", "Cleveland", 0, 3, 7, 0, 10, 48, 47, 180, 1.5),
  @(1211, 114, "H", "Cincinnati", 0, 13, 7, 3, 23, 3.5, 4, -210, 23.5),
  @(1211, 115, "V", "Jacksonville", 7, 13, 13, 3, 36, 43, 41.5, 160, 21.5),
  @(1211, 116, "H", "Tennessee", 14, 0, 0, 8, 22, 4, 3, -180, 3),
  @(1211, 117, "V", "Houston", 10, 10, 3, 0, 23, 47.5, 44.5, 1100, 23),
  @(1211, 118, "H", "Dallas", 7, 10, 0, 10, 27, 14, 17, -2500, 10),
  @(1211, 119, "V", "Miami", 0, 7, 7, 3, 17, 51.5, 3, -170, 4),
  @(1211, 120, "H", "LAChargers", 0, 17, 3, 3, 23, 1, 55, 150, 27),
  @(1211, 121, "V", "Carolina", 10, 10, 0, 10, 30, 43.5, 44, 170, 22),
  @(1211, 122, "H", "Seattle", 0, 14, 3, 7, 24, 7, 4, -190, 4.5),
  @(1211, 123, "V", "TampaBay", 0, 0, 7, 0, 7, 41.5, 39, 165, 3),
  @(1211, 124, "H", "SanFrancisco", 7, 21, 7, 0, 35, 5, 3.5, -185, 17),
  @(1211, 125, "V", "KansasCity", 6, 21, 7, 0, 34, 6, 9.5, -420, 3),
  @(1211, 126, "H", "Denver", 0, 14, 7, 7, 28, 42.5, 44, 350, 22),
)

$startRow = 390
$afterRow = $startRow + $data.Count

# The row right after the new data block (previously just another blank,
# specially-formatted placeholder row, like the ones the new rows replace)
# needs to collapse back to a plain, implicit default row with no explicit
# formatting of its own. Doing a delete+insert round trip *before* any other
# edits forces the row metadata to be recomputed so it fully disappears from
# the saved XML, instead of copying formatting from its neighbors.
$ws.Rows.Item($afterRow).Delete()
$ws.Rows.Item($afterRow).Insert()

# Clear any leftover row-level formatting (customFormat/s) on the target rows
# before writing values, so the new rows end up without a row-level style
# override (matching freshly-authored data rows), same as the rows above them.
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $ws.Rows.Item($r).ClearFormats()
}

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

Write-Output "done"